# Update "想去人数" (interested-count) figures in the F column across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets.
# 本地生活 (sheet3) has no changes in this revision.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 5508
$ws1.Range("F9").Value  = 9
$ws1.Range("F11").Value = 816
$ws1.Range("F12").Value = 24
$ws1.Range("F13").Value = 6614
$ws1.Range("F14").Value = 43
$ws1.Range("F17").Value = 5104
$ws1.Range("F18").Value = 115
$ws1.Range("F20").Value = 4226
$ws1.Range("F22").Value = 4175
$ws1.Range("F23").Value = 214
$ws1.Range("F26").Value = 290
$ws1.Range("F27").Value = 266
$ws1.Range("F32").Value = 63
$ws1.Range("F33").Value = 7564
$ws1.Range("F35").Value = 1275
$ws1.Range("F36").Value = 629
$ws1.Range("F38").Value = 981
$ws1.Range("F40").Value = 1506
$ws1.Range("F42").Value = 847
$ws1.Range("F44").Value = 3725
$ws1.Range("F45").Value = 335
$ws1.Range("F48").Value = 818
$ws1.Range("F49").Value = 1045

# --- 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 101
$ws2.Range("F17").Value = 3
$ws2.Range("F18").Value = 70

# --- 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 5508
$ws4.Range("F8").Value  = 101
$ws4.Range("F11").Value = 9
$ws4.Range("F14").Value = 816
$ws4.Range("F15").Value = 24
$ws4.Range("F16").Value = 6614
$ws4.Range("F17").Value = 43
$ws4.Range("F20").Value = 5104
$ws4.Range("F21").Value = 115
$ws4.Range("F23").Value = 4226
$ws4.Range("F24").Value = 4175
$ws4.Range("F25").Value = 214
$ws4.Range("F27").Value = 290
$ws4.Range("F28").Value = 266
$ws4.Range("F32").Value = 7565
$ws4.Range("F34").Value = 1275
$ws4.Range("F35").Value = 629
$ws4.Range("F37").Value = 981
$ws4.Range("F39").Value = 1506
$ws4.Range("F41").Value = 847
$ws4.Range("F43").Value = 3725
$ws4.Range("F44").Value = 335
$ws4.Range("F47").Value = 818
$ws4.Range("F48").Value = 1045
